$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 18; $r++) {
    $ws.Range("AO$r").Value = 321688.56653439248
}

$excel.Calculate()
